$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F241").Value = 89610
$ws.Range("G273").Value = 1291
$ws.Range("F279").Value = 43435
$ws.Range("G279").Value = 3116
$ws.Range("F281").Value = 45603
$ws.Range("F282").Value = 47426
$ws.Range("G282").Value = 2840
$ws.Range("F283").Value = 16932
$ws.Range("G283").Value = 1003
$ws.Range("F285").Value = 42163
$ws.Range("G285").Value = 3431
$ws.Range("F286").Value = 55435
$ws.Range("G286").Value = 4298
$ws.Range("F287").Value = 59001
$ws.Range("G287").Value = 3740
$ws.Range("F288").Value = 58669
$ws.Range("G288").Value = 4047
$ws.Range("F289").Value = 63132
$ws.Range("G289").Value = 3658
$ws.Range("F291").Value = 14954
$ws.Range("G291").Value = 486
$ws.Range("F292").Value = 82431
$ws.Range("G292").Value = 7273
$ws.Range("F293").Value = 82339
$ws.Range("G293").Value = 5762
$ws.Range("F294").Value = 93683
$ws.Range("G294").Value = 4930
$ws.Range("F299").Value = 65527
$ws.Range("G299").Value = 6887
$ws.Range("F300").Value = 72345
$ws.Range("G300").Value = 7062
$ws.Range("F301").Value = 71778
$ws.Range("G301").Value = 5645
$ws.Range("F302").Value = 77386
$ws.Range("G302").Value = 5605
$ws.Range("F306").Value = 73301
$ws.Range("G306").Value = 7346
$ws.Range("F307").Value = 77479
$ws.Range("G307").Value = 6578
$ws.Range("F308").Value = 15623
$ws.Range("G308").Value = 1084
$ws.Range("F309").Value = 76770
$ws.Range("G309").Value = 5393
$ws.Range("F310").Value = 78001
$ws.Range("G310").Value = 4054
$ws.Range("F311").Value = 62148
$ws.Range("G311").Value = 1971
$ws.Range("F312").Value = 28009
$ws.Range("G312").Value = 926
$ws.Range("F313").Value = 73038
$ws.Range("G313").Value = 3283
$ws.Range("F314").Value = 64582
$ws.Range("G314").Value = 3194
$ws.Range("F315").Value = 56898
$ws.Range("G315").Value = 2674
$ws.Range("F316").Value = 50387
$ws.Range("G316").Value = 2274
$ws.Range("F317").Value = 63349
$ws.Range("G317").Value = 2180
$ws.Range("F318").Value = 49865
$ws.Range("G318").Value = 1198
$ws.Range("F319").Value = 41308
$ws.Range("F320").Value = 70646
$ws.Range("G320").Value = 3208
$ws.Range("F321").Value = 94034
$ws.Range("G321").Value = 2859
$ws.Range("F322").Value = 109484
$ws.Range("G322").Value = 2355
$ws.Range("F323").Value = 216313
$ws.Range("G323").Value = 3218
$ws.Range("F324").Value = 237729
$ws.Range("G324").Value = 2750
$ws.Range("F325").Value = 756934
$ws.Range("G325").Value = 6359
$ws.Range("F326").Value = 437574
$ws.Range("G326").Value = 3917
$ws.Range("F327").Value = 235098
$ws.Range("G327").Value = 2852
$ws.Range("F328").Value = 179871
$ws.Range("G328").Value = 2650
$ws.Range("F329").Value = 83325
$ws.Range("G329").Value = 1759
$ws.Range("F330").Value = 71915
$ws.Range("G330").Value = 2046
$ws.Range("F331").Value = 151925
$ws.Range("G331").Value = 2656
$ws.Range("F332").Value = 444402
$ws.Range("G332").Value = 4426
$ws.Range("F333").Value = 270559
$ws.Range("G333").Value = 2903
$ws.Range("F334").Value = 203690
$ws.Range("G334").Value = 3411
$ws.Range("F335").Value = 131373
$ws.Range("G335").Value = 2969
$ws.Range("F336").Value = 102904
$ws.Range("G336").Value = 3277
$ws.Range("F337").Value = 104782
$ws.Range("G337").Value = 2969
$ws.Range("F338").Value = 222374
$ws.Range("G338").Value = 3111
$ws.Range("F339").Value = 645272
$ws.Range("G339").Value = 5503
$ws.Range("F340").Value = 380511
$ws.Range("G340").Value = 3275
$ws.Range("F341").Value = 292024
$ws.Range("G341").Value = 3601
$ws.Range("F342").Value = 174508
$ws.Range("G342").Value = 2953
$ws.Range("F343").Value = 127095
$ws.Range("G343").Value = 2829
$ws.Range("F344").Value = 131777
$ws.Range("G344").Value = 2432
$ws.Range("F345").Value = 279564
$ws.Range("G345").Value = 3199
$ws.Range("F346").Value = 647469
$ws.Range("G346").Value = 4621
$ws.Range("F347").Value = 327484
$ws.Range("G347").Value = 2769
$ws.Range("F348").Value = 224684
$ws.Range("G348").Value = 3092
$ws.Range("F349").Value = 152143
$ws.Range("G349").Value = 2619
$ws.Range("F350").Value = 121363
$ws.Range("G350").Value = 2634
$ws.Range("F351").Value = 140455
$ws.Range("G351").Value = 2586
$ws.Range("F352").Value = 288495
$ws.Range("G352").Value = 3782
$ws.Range("F353").Value = 667670
$ws.Range("G353").Value = 4928
$ws.Range("F354").Value = 283325
$ws.Range("G354").Value = 2821
$ws.Range("F355").Value = 206358
$ws.Range("G355").Value = 3148
$ws.Range("F356").Value = 132614
$ws.Range("G356").Value = 2323
